$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all data rows (2-234)
# from serial 46061 (2026-02-08) to 46062 (2026-02-09).
$ws.Range("C2:C234").Value = 46062
